$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.102882404376526
$ws.Cells.Item(2, 4).Value = 1.111984274116987
$ws.Cells.Item(2, 5).Value = 1.104158919159907
$ws.Cells.Item(2, 6).Value = 1.117494588057065
$ws.Cells.Item(2, 9).Value = 1.063851322556931
$ws.Cells.Item(2, 10).Value = 1.107653188708778
$ws.Cells.Item(2, 11).Value = 1.114577768308342
$ws.Cells.Item(2, 12).Value = 1.106771784569194
$ws.Cells.Item(2, 13).Value = 1.120074616826322
$ws.Cells.Item(2, 14).Value = 1.109226183286341

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.104680209325298
$ws.Cells.Item(3, 4).Value = 1.113706973285099
$ws.Cells.Item(3, 5).Value = 1.105784844255862
$ws.Cells.Item(3, 6).Value = 1.119211531459812
$ws.Cells.Item(3, 9).Value = 1.064364894093118
$ws.Cells.Item(3, 10).Value = 1.109116791844593
$ws.Cells.Item(3, 11).Value = 1.116122400155308
$ws.Cells.Item(3, 12).Value = 1.108218446758596
$ws.Cells.Item(3, 13).Value = 1.121614491499204
$ws.Cells.Item(3, 14).Value = 1.110691864906486

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.105840661371052
$ws.Cells.Item(4, 4).Value = 1.114819074885089
$ws.Cells.Item(4, 5).Value = 1.106834081408451
$ws.Cells.Item(4, 6).Value = 1.120319892919239
$ws.Cells.Item(4, 9).Value = 1.064694348096832
$ws.Cells.Item(4, 10).Value = 1.110060605276339
$ws.Cells.Item(4, 11).Value = 1.117118799421269
$ws.Cells.Item(4, 12).Value = 1.109151193053072
$ws.Cells.Item(4, 13).Value = 1.122607793486566
$ws.Cells.Item(4, 14).Value = 1.111637018661562

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.106327848991577
$ws.Cells.Item(5, 4).Value = 1.115285994063125
$ws.Cells.Item(5, 5).Value = 1.107274513299948
$ws.Cells.Item(5, 6).Value = 1.120785235911389
$ws.Cells.Item(5, 9).Value = 1.064832170834366
$ws.Cells.Item(5, 10).Value = 1.110456622595924
$ws.Cells.Item(5, 11).Value = 1.117536961340993
$ws.Cells.Item(5, 12).Value = 1.109542532683325
$ws.Cells.Item(5, 13).Value = 1.123024648716195
$ws.Cells.Item(5, 14).Value = 1.112033598371165

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.1064096113024
$ws.Cells.Item(6, 4).Value = 1.115364356582072
$ws.Cells.Item(6, 5).Value = 1.107348425011855
$ws.Cells.Item(6, 6).Value = 1.120863333555597
$ws.Cells.Item(6, 9).Value = 1.064855272170074
$ws.Cells.Item(6, 10).Value = 1.110523071340546
$ws.Cells.Item(6, 11).Value = 1.117607130474261
$ws.Cells.Item(6, 12).Value = 1.109608194557146
$ws.Cells.Item(6, 13).Value = 1.123094598180215
$ws.Cells.Item(6, 14).Value = 1.112100141480626

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.105847173793381
$ws.Cells.Item(7, 4).Value = 1.114825316253
$ws.Cells.Item(7, 5).Value = 1.106839969084588
$ws.Cells.Item(7, 6).Value = 1.120326113241201
$ws.Cells.Item(7, 9).Value = 1.064696192354642
$ws.Cells.Item(7, 10).Value = 1.110065899855159
$ws.Cells.Item(7, 11).Value = 1.117124389749298
$ws.Cells.Item(7, 12).Value = 1.109156425226704
$ws.Cells.Item(7, 13).Value = 1.122613366372629
$ws.Cells.Item(7, 14).Value = 1.111642320759292

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.103490579969065
$ws.Cells.Item(8, 4).Value = 1.1125670152263
$ws.Cells.Item(8, 5).Value = 1.104709005546414
$ws.Cells.Item(8, 6).Value = 1.118075387212083
$ws.Cells.Item(8, 9).Value = 1.064025482040965
$ws.Cells.Item(8, 10).Value = 1.108148497930307
$ws.Cells.Item(8, 11).Value = 1.115100429771191
$ws.Cells.Item(8, 12).Value = 1.107261390376755
$ws.Cells.Item(8, 13).Value = 1.12059567451051
$ws.Cells.Item(8, 14).Value = 1.109722195903778

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.099315455066291
$ws.Cells.Item(9, 4).Value = 1.108567052327218
$ws.Cells.Item(9, 5).Value = 1.100931567405441
$ws.Cells.Item(9, 6).Value = 1.114088654069486
$ws.Cells.Item(9, 9).Value = 1.062821455749108
$ws.Cells.Item(9, 10).Value = 1.104744433494473
$ws.Cells.Item(9, 11).Value = 1.111509762239837
$ws.Cells.Item(9, 12).Value = 1.103895928444343
$ws.Cells.Item(9, 13).Value = 1.117015907496045
$ws.Cells.Item(9, 14).Value = 1.106313297305993

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.096515933401987
$ws.Cells.Item(10, 4).Value = 1.10588570258103
$ws.Cells.Item(10, 5).Value = 1.098397341509095
$ws.Cells.Item(10, 6).Value = 1.11141604807559
$ws.Cells.Item(10, 9).Value = 1.062003547060509
$ws.Cells.Item(10, 10).Value = 1.102457199235229
$ws.Cells.Item(10, 11).Value = 1.109098881597639
$ws.Cells.Item(10, 12).Value = 1.1016338766811
$ws.Cells.Item(10, 13).Value = 1.114612198499366
$ws.Cells.Item(10, 14).Value = 1.104022814911751

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.095299671916041
$ws.Cells.Item(11, 4).Value = 1.104720960283738
$ws.Cells.Item(11, 5).Value = 1.097296020264916
$ws.Cells.Item(11, 6).Value = 1.110255074662407
$ws.Cells.Item(11, 9).Value = 1.061645697770973
$ws.Cells.Item(11, 10).Value = 1.101462384715778
$ws.Cells.Item(11, 11).Value = 1.108050699787543
$ws.Cells.Item(11, 12).Value = 1.100649833918955
$ws.Cells.Item(11, 13).Value = 1.113567099914101
$ws.Cells.Item(11, 14).Value = 1.103026587641576

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.094847270899307
$ws.Cells.Item(12, 4).Value = 1.104287749917477
$ws.Cells.Item(12, 5).Value = 1.096886325037579
$ws.Cells.Item(12, 6).Value = 1.10982326170905
$ws.Cells.Item(12, 9).Value = 1.061512215794582
$ws.Cells.Item(12, 10).Value = 1.101092185491654
$ws.Cells.Item(12, 11).Value = 1.107660702937947
$ws.Cells.Item(12, 12).Value = 1.100283615786325
$ws.Cells.Item(12, 13).Value = 1.11317824502318
$ws.Cells.Item(12, 14).Value = 1.102655862692091

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.094944341146938
$ws.Cells.Item(13, 4).Value = 1.104380701230874
$ws.Cells.Item(13, 5).Value = 1.096974234178501
$ws.Cells.Item(13, 6).Value = 1.109915913385489
$ws.Cells.Item(13, 9).Value = 1.061540873615332
$ws.Cells.Item(13, 10).Value = 1.101171625544487
$ws.Cells.Item(13, 11).Value = 1.107744388509677
$ws.Cells.Item(13, 12).Value = 1.100362202789355
$ws.Cells.Item(13, 13).Value = 1.113261685789068
$ws.Cells.Item(13, 14).Value = 1.102735415558911

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.095262289214421
$ws.Cells.Item(14, 4).Value = 1.104685162774009
$ws.Cells.Item(14, 5).Value = 1.097262167383238
$ws.Cells.Item(14, 6).Value = 1.110219392714239
$ws.Cells.Item(14, 9).Value = 1.061634675587675
$ws.Cells.Item(14, 10).Value = 1.101431797912433
$ws.Cells.Item(14, 11).Value = 1.108018475992058
$ws.Cells.Item(14, 12).Value = 1.100619576602951
$ws.Cells.Item(14, 13).Value = 1.113534970580128
$ws.Cells.Item(14, 14).Value = 1.102995957401462

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.095458103814997
$ws.Cells.Item(15, 4).Value = 1.104872675101241
$ws.Cells.Item(15, 5).Value = 1.097439490469963
$ws.Cells.Item(15, 6).Value = 1.110406299535122
$ws.Cells.Item(15, 9).Value = 1.061692395568525
$ws.Cells.Item(15, 10).Value = 1.101592008007058
$ws.Cells.Item(15, 11).Value = 1.108187262992461
$ws.Cells.Item(15, 12).Value = 1.100778059759571
$ws.Cells.Item(15, 13).Value = 1.113703262581363
$ws.Cells.Item(15, 14).Value = 1.103156395012797

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.096596565734537
$ws.Cells.Item(16, 4).Value = 1.105962923260168
$ws.Cells.Item(16, 5).Value = 1.098470347188031
$ws.Cells.Item(16, 6).Value = 1.111493018269337
$ws.Cells.Item(16, 9).Value = 1.062027218052383
$ws.Cells.Item(16, 10).Value = 1.102523127189886
$ws.Cells.Item(16, 11).Value = 1.109168354950109
$ws.Cells.Item(16, 12).Value = 1.101699086949024
$ws.Cells.Item(16, 13).Value = 1.114681466765069
$ws.Cells.Item(16, 14).Value = 1.104088836491666

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.097309595235241
$ws.Cells.Item(17, 4).Value = 1.106645804329101
$ws.Cells.Item(17, 5).Value = 1.099115897734616
$ws.Cells.Item(17, 6).Value = 1.112173680861397
$ws.Cells.Item(17, 9).Value = 1.062236251005912
$ws.Cells.Item(17, 10).Value = 1.103105997385757
$ws.Cells.Item(17, 11).Value = 1.109782617526683
$ws.Cells.Item(17, 12).Value = 1.102275591161421
$ws.Cells.Item(17, 13).Value = 1.115293912015476
$ws.Cells.Item(17, 14).Value = 1.104672534430071

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.097725103785109
$ws.Cells.Item(18, 4).Value = 1.107043761470199
$ws.Cells.Item(18, 5).Value = 1.099492053125881
$ws.Cells.Item(18, 6).Value = 1.112570342336457
$ws.Cells.Item(18, 9).Value = 1.06235782085934
$ws.Cells.Item(18, 10).Value = 1.103445549639465
$ws.Cells.Item(18, 11).Value = 1.110140497096682
$ws.Cells.Item(18, 12).Value = 1.102611417475217
$ws.Cells.Item(18, 13).Value = 1.115650729430518
$ws.Cells.Item(18, 14).Value = 1.105012568886926

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.09786671596536
$ws.Cells.Item(19, 4).Value = 1.107179394736457
$ws.Cells.Item(19, 5).Value = 1.099620247881911
$ws.Cells.Item(19, 6).Value = 1.112705533529447
$ws.Cells.Item(19, 9).Value = 1.062399212929421
$ws.Cells.Item(19, 10).Value = 1.103561256386927
$ws.Cells.Item(19, 11).Value = 1.110262455829451
$ws.Cells.Item(19, 12).Value = 1.102725851611309
$ws.Cells.Item(19, 13).Value = 1.115772325631871
$ws.Cells.Item(19, 14).Value = 1.105128439951242

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.097233134341214
$ws.Cells.Item(20, 4).Value = 1.106572574583542
$ws.Cells.Item(20, 5).Value = 1.099046676040498
$ws.Cells.Item(20, 6).Value = 1.112100689309244
$ws.Cells.Item(20, 9).Value = 1.062213860569661
$ws.Cells.Item(20, 10).Value = 1.103043505113811
$ws.Cells.Item(20, 11).Value = 1.109716755428846
$ws.Cells.Item(20, 12).Value = 1.102213783220953
$ws.Cells.Item(20, 13).Value = 1.115228245120544
$ws.Cells.Item(20, 14).Value = 1.10460995341193

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.095168678903561
$ws.Cells.Item(21, 4).Value = 1.104595522401639
$ws.Cells.Item(21, 5).Value = 1.097177395357319
$ws.Cells.Item(21, 6).Value = 1.11013004164849
$ws.Cells.Item(21, 9).Value = 1.061607068778837
$ws.Cells.Item(21, 10).Value = 1.101355202555667
$ws.Cells.Item(21, 11).Value = 1.107937782285154
$ws.Cells.Item(21, 12).Value = 1.100543805904277
$ws.Cells.Item(21, 13).Value = 1.113454513336268
$ws.Cells.Item(21, 14).Value = 1.102919253270503

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.09386703470017
$ws.Cells.Item(22, 4).Value = 1.103349146053112
$ws.Cells.Item(22, 5).Value = 1.095998534947771
$ws.Cells.Item(22, 6).Value = 1.10888767754047
$ws.Cells.Item(22, 9).Value = 1.061222306694092
$ws.Cells.Item(22, 10).Value = 1.100289752815943
$ws.Cells.Item(22, 11).Value = 1.106815470817257
$ws.Cells.Item(22, 12).Value = 1.099489762089898
$ws.Cells.Item(22, 13).Value = 1.112335478313955
$ws.Cells.Item(22, 14).Value = 1.101852290469941

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.094557412219596
$ws.Cells.Item(23, 4).Value = 1.104010194755246
$ws.Cells.Item(23, 5).Value = 1.096623815163723
$ws.Cells.Item(23, 6).Value = 1.109546600632647
$ws.Cells.Item(23, 9).Value = 1.061426586426055
$ws.Cells.Item(23, 10).Value = 1.100854947295573
$ws.Cells.Item(23, 11).Value = 1.107410795073189
$ws.Cells.Item(23, 12).Value = 1.100048921097177
$ws.Cells.Item(23, 13).Value = 1.112929067452325
$ws.Cells.Item(23, 14).Value = 1.102418287590559

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.097267684931239
$ws.Cells.Item(24, 4).Value = 1.10660566505271
$ws.Cells.Item(24, 5).Value = 1.099077955528672
$ws.Cells.Item(24, 6).Value = 1.112133672154215
$ws.Cells.Item(24, 9).Value = 1.062223978942513
$ws.Cells.Item(24, 10).Value = 1.103071743997871
$ws.Cells.Item(24, 11).Value = 1.109746516940659
$ws.Cells.Item(24, 12).Value = 1.102241712924518
$ws.Cells.Item(24, 13).Value = 1.115257918435177
$ws.Cells.Item(24, 14).Value = 1.104638232398446

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.100397591559069
$ws.Cells.Item(25, 4).Value = 1.109603664960432
$ws.Cells.Item(25, 5).Value = 1.101910869980009
$ws.Cells.Item(25, 6).Value = 1.115121860164508
$ws.Cells.Item(25, 9).Value = 1.063135384194622
$ws.Cells.Item(25, 10).Value = 1.105627553486291
$ws.Cells.Item(25, 11).Value = 1.112440988331694
$ws.Cells.Item(25, 12).Value = 1.104769165633658
$ws.Cells.Item(25, 13).Value = 1.117944332496756
$ws.Cells.Item(25, 14).Value = 1.107197671429494
